$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 99. This pushes the existing row 99
# (and everything below it, down through the old row 151) down by one row,
# so the old row 151 becomes the new row 152 - matching the diff.
$ws.Rows.Item(99).Insert()

# Populate the newly-inserted row 99 with the new weekly record.
$ws.Range("A99").Value = 4
$ws.Range("B99").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C99").Value = "Los Lagos"
$ws.Range("D99").Value = 44466
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E99").Value = 10
$ws.Range("F99").Value = 100112044
$ws.Range("G99").Value = "Perejil"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 90
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 4500
$ws.Range("M99").Value = 4500
$ws.Range("N99").Value = '$/docena de atados (3 kilos)'
$ws.Range("O99").Value = "Región Metropolitana"
$ws.Range("P99").Value = 1500
$ws.Range("Q99").Value = 3
$ws.Range("R99").Value = "Hortaliza"
